# chore: adapt column header formatting to respective input file names (#7)
#
# Renames the "_old" / "_new" column-header suffixes to the respective
# input-file-version suffixes ("_FV2410" / "_FV2504"), turns the header
# row + data range into a native Excel Table (ListObject), and freezes
# the header row so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row (row 1, columns A:U) -----------------------------
# Columns A-J carried the "_old" suffix -> "_FV2410"
# Column  K is the untouched "diff" header
# Columns L-U carried the "_new" suffix -> "_FV2504"
$oldHeaders = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)
$newHeaders = @(
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)

for ($i = 0; $i -lt $oldHeaders.Length; $i++) {
    # columns 1..10 => A..J
    $ws.Cells.Item(1, $i + 1).Value = $oldHeaders[$i]
}
for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    # columns 12..21 => L..U
    $ws.Cells.Item(1, $i + 12).Value = $newHeaders[$i]
}

# --- 2. Turn A1:U56 into a native Excel table ------------------------------
$range = $ws.Range("A1:U56")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $range, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# --- 3. Freeze the header row ----------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
